$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 32 (pushes existing rows 32:121 down to 33:122,
# which matches every other data row shifting down by one and the final
# existing row 121 duplicating into the new row 122).
$ws.Rows("32:32").Insert()

# Populate the newly inserted row 32 with its data (same constant columns as
# every other data row in this sheet, plus the row-specific values).
$ws.Cells.Item(32, 1).Value2  = 1
$ws.Cells.Item(32, 2).Value2  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(32, 3).Value2  = "Arica y Parinacota"
$ws.Cells.Item(32, 4).Value2  = 45259
$ws.Cells.Item(32, 5).Value2  = 15
$ws.Cells.Item(32, 6).Value2  = 100112012
$ws.Cells.Item(32, 7).Value2  = "Espinaca"
$ws.Cells.Item(32, 8).Value2  = "Sin especificar"
$ws.Cells.Item(32, 9).Value2  = "Segunda"
$ws.Cells.Item(32, 10).Value2 = 400
$ws.Cells.Item(32, 11).Value2 = 1500
$ws.Cells.Item(32, 12).Value2 = 2000
$ws.Cells.Item(32, 13).Value2 = 1688
$ws.Cells.Item(32, 14).Value2 = "$/atado 2,5 a 3 kilos"
$ws.Cells.Item(32, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(32, 16).Value2 = 563
$ws.Cells.Item(32, 17).Value2 = 3
$ws.Cells.Item(32, 18).Value2 = "Hortaliza"

# Make sure the date cell keeps the date/time number format used by every
# other row's column D (style index 2 in the original workbook).
$ws.Cells.Item(32, 4).NumberFormat = $ws.Cells.Item(33, 4).NumberFormat
